$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# --- Row 7 : update interest/capital split for period 3 ---
$ws1.Range("D7").Value = 16006.36
$ws1.Range("E7").Value = 1993.64

# --- Row 8 : update interest/capital split for period 4 ---
$ws1.Range("D8").Value = 6406.52
$ws1.Range("E8").Value = 1593.48

# --- Row 9 : period 5 now has a payment recorded ---
$ws1.Range("C9").Formula = "=C8-D9"
$ws1.Range("D9").Value = 7566.68
$ws1.Range("E9").Value = 1433.32
$ws1.Range("G9").Formula = "=(`$D`$2-SUM(`$D`$5:D8))*30/100/12"

# --- Row 10 : period 6 now has a payment recorded ---
$ws1.Range("C10").Formula = "=C9-D10"
$ws1.Range("D10").Value = 13755.85
$ws1.Range("E10").Value = 1244.15
$ws1.Range("F10").Formula = "=D10+E10"
$ws1.Range("G10").Formula = "=(`$D`$2-SUM(`$D`$5:D9))*30/100/12"

# --- Row 11 : period 7 now has a payment recorded ---
$ws1.Range("C11").Formula = "=C10-D11"
$ws1.Range("D11").Value = 3099.74
$ws1.Range("E11").Value = 900.26
$ws1.Range("G11").Formula = "=(`$D`$2-SUM(`$D`$5:D10))*30/100/12"

# --- Row 12 : period 8 now has a payment recorded ---
$ws1.Range("C12").Formula = "=C11-D12"
$ws1.Range("D12").Value = 4000
$ws1.Range("G12").Formula = "=(`$D`$2-SUM(`$D`$5:D13))*30/100/12"

# --- Row 13 : period 9 now has a payment recorded ---
$ws1.Range("C13").Formula = "=C12-D13"
$ws1.Range("D13").Value = 6000

# --- Sheet / selection bookkeeping: Hoja1 becomes the active tab ---
# (Hoja2's own selection, H12, is untouched by this edit - only its
#  tabSelected flag changes as a side effect of Hoja1 becoming active.)
$ws1.Activate() | Out-Null
$ws1.Range("G13").Select() | Out-Null
